# Corrected pH inputs after talk with Irith
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry pH")
$ws.Activate()

$ws.Range("A2").Value = 8.3
$ws.Range("A3").Value = 8.1
$ws.Range("A4").Value = 7.9
$ws.Range("A5").Value = 7.7
$ws.Range("A6").Value = 7.5
$ws.Range("A7").Value = 7.3
$ws.Range("A8").Value = 7.1

$ws.Range("N4").Select()
